$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DMD")
$ws.Activate()

# Clear the "TimeSlice" helper table (rows 26-40, columns B-D) entirely
# (values, formulas and formatting) as it was removed from the model.
$ws.Range("B26:E40").Clear()

$ws.Range("F33").Select()
